$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26; this shifts the existing rows 26-32 down to 27-33,
# preserving all of their data (matching the target diff).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44795
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112045
$ws.Cells.Item(26, 7).Value = "Zapallo"
$ws.Cells.Item(26, 8).Value = "Camote"
$ws.Cells.Item(26, 9).Value = "1a nueva(o)"
$ws.Cells.Item(26, 10).Value = 1000
$ws.Cells.Item(26, 11).Value = 1000
$ws.Cells.Item(26, 12).Value = 1100
$ws.Cells.Item(26, 13).Value = 1050
$ws.Cells.Item(26, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 1050
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
